# Venv를 통한 가상 환경.pptx - edit script
#
# Target change (per commit diff):
#   1. Delete the trailing slide (slide 8 / sldId 296 / rId12 -> ppt/slides/slide8.xml),
#      a simple "thank you / contact info" closing slide. This also drops the now
#      unused slide relationship (rId12) from presentation.xml.rels, which in turn
#      shifts the notesMaster/handoutMaster relationship ids down by one
#      (rId13->rId12, rId14->rId13) once PowerPoint renumbers them on save.
#   2. Refresh the cached "datetime1" field text (an auto date field) on both the
#      Notes Master and the Handout Master from 2023-03-10 to 2023-03-14.

$p = $ppt.ActivePresentation

# --- 1. Remove the last slide -------------------------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()

# --- 2. Update the cached date field text on Notes Master & Handout Master ---
$newDate = "2023-03-14"

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $shape = $nm.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "2023-03-10") {
            $tr.Text = $newDate
        }
    }
}

$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $shape = $hm.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "2023-03-10") {
            $tr.Text = $newDate
        }
    }
}
